{"js": "// The document has a 6x5 table of \"A\u00f7B=C, R\" cells. This edit swaps each\n// cell's three-digit-number \u00f7 one-digit-number problem/answer for a new\n// one. Every \"old\" text below occurs exactly once in the document, so an\n// exact (case-sensitive) search unambiguously finds the run to replace,\n// and insertText(..., \"Replace\") swaps its text while leaving the run's\n// formatting (font/size) untouched.\nconst replacements = [\n  [\"361\u00f75=72, 1\", \"750\u00f73=250, 0\"],\n  [\"391\u00f74=97, 3\", \"815\u00f76=135, 5\"],\n  [\"791\u00f79=87, 8\", \"952\u00f73=317, 1\"],\n  [\"261\u00f76=43, 3\", \"611\u00f76=101, 5\"],\n  [\"408\u00f77=58, 2\", \"184\u00f72=92, 0\"],\n  [\"263\u00f73=87, 2\", \"369\u00f74=92, 1\"],\n  [\"344\u00f79=38, 2\", \"814\u00f75=162, 4\"],\n  [\"956\u00f78=119, 4\", \"542\u00f77=77, 3\"],\n  [\"183\u00f72=91, 1\", \"745\u00f75=149, 0\"],\n  [\"697\u00f75=139, 2\", \"260\u00f74=65, 0\"],\n  [\"660\u00f74=165, 0\", \"771\u00f74=192, 3\"],\n  [\"807\u00f76=134, 3\", \"682\u00f72=341, 0\"],\n  [\"160\u00f75=32, 0\", \"140\u00f75=28, 0\"],\n  [\"253\u00f77=36, 1\", \"553\u00f77=79, 0\"],\n  [\"175\u00f79=19, 4\", \"322\u00f77=46, 0\"],\n  [\"315\u00f74=78, 3\", \"111\u00f79=12, 3\"],\n  [\"502\u00f78=62, 6\", \"715\u00f73=238, 1\"],\n  [\"923\u00f76=153, 5\", \"388\u00f72=194, 0\"],\n  [\"519\u00f75=103, 4\", \"267\u00f79=29, 6\"],\n  [\"921\u00f78=115, 1\", \"379\u00f76=63, 1\"],\n  [\"584\u00f76=97, 2\", \"861\u00f72=430, 1\"],\n  [\"454\u00f72=227, 0\", \"991\u00f76=165, 1\"],\n  [\"121\u00f79=13, 4\", \"991\u00f79=110, 1\"],\n  [\"120\u00f78=15, 0\", \"186\u00f72=93, 0\"],\n  [\"464\u00f75=92, 4\", \"387\u00f73=129, 0\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items.forEach((item) => {\n    item.insertText(newText, Word.InsertLocation.replace);\n  });\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-division expression with its new value.\n# Every source string is unique in the document, so Find.Execute with\n# MatchCase + Wrap=wdFindStop safely targets exactly one run each time.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"361\u00f75=72, 1\", \"750\u00f73=250, 0\"),\n  @(\"391\u00f74=97, 3\", \"815\u00f76=135, 5\"),\n  @(\"791\u00f79=87, 8\", \"952\u00f73=317, 1\"),\n  @(\"261\u00f76=43, 3\", \"611\u00f76=101, 5\"),\n  @(\"408\u00f77=58, 2\", \"184\u00f72=92, 0\"),\n  @(\"263\u00f73=87, 2\", \"369\u00f74=92, 1\"),\n  @(\"344\u00f79=38, 2\", \"814\u00f75=162, 4\"),\n  @(\"956\u00f78=119, 4\", \"542\u00f77=77, 3\"),\n  @(\"183\u00f72=91, 1\", \"745\u00f75=149, 0\"),\n  @(\"697\u00f75=139, 2\", \"260\u00f74=65, 0\"),\n  @(\"660\u00f74=165, 0\", \"771\u00f74=192, 3\"),\n  @(\"807\u00f76=134, 3\", \"682\u00f72=341, 0\"),\n  @(\"160\u00f75=32, 0\", \"140\u00f75=28, 0\"),\n  @(\"253\u00f77=36, 1\", \"553\u00f77=79, 0\"),\n  @(\"175\u00f79=19, 4\", \"322\u00f77=46, 0\"),\n  @(\"315\u00f74=78, 3\", \"111\u00f79=12, 3\"),\n  @(\"502\u00f78=62, 6\", \"715\u00f73=238, 1\"),\n  @(\"923\u00f76=153, 5\", \"388\u00f72=194, 0\"),\n  @(\"519\u00f75=103, 4\", \"267\u00f79=29, 6\"),\n  @(\"921\u00f78=115, 1\", \"379\u00f76=63, 1\"),\n  @(\"584\u00f76=97, 2\", \"861\u00f72=430, 1\"),\n  @(\"454\u00f72=227, 0\", \"991\u00f76=165, 1\"),\n  @(\"121\u00f79=13, 4\", \"991\u00f79=110, 1\"),\n  @(\"120\u00f78=15, 0\", \"186\u00f72=93, 0\"),\n  @(\"464\u00f75=92, 4\", \"387\u00f73=129, 0\")\n)\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n#   ReplaceWith, Replace) -- Wrap=1 is wdFindStop, Replace=2 is wdReplaceAll.\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}\n\n"}
